$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1541.4615
$ws.Range("I40").Value = 1173.3334
$ws.Range("K40").Value = 1173.3334
$ws.Range("M40").Value = -998.3334
$ws.Range("H43").Value = 2133.3333
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 2133.3333
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 2133.3333
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -2271.3333
$ws.Range("H64").Value = 3804.842
$ws.Range("I64").Value = 3207.6667
$ws.Range("K64").Value = 3207.6667
$ws.Range("M64").Value = -2959.6667
$ws.Range("H67").Value = 3804.842
$ws.Range("I67").Value = 3207.6667
$ws.Range("K67").Value = 3207.6667
$ws.Range("M67").Value = -2349.6667
$ws.Range("H76").Value = 3089472.8
$ws.Range("I76").Value = 3150.3845
$ws.Range("K76").Value = 3150.3845
$ws.Range("M76").Value = -2835.3845
$ws.Range("H79").Value = 3089472.8
$ws.Range("I79").Value = 3150.3845
$ws.Range("K79").Value = 3150.3845
$ws.Range("M79").Value = -2058.3845
$ws.Range("H86").Value = 6936.0527
$ws.Range("I86").Value = 1685.4445
$ws.Range("J86").Value = 11661.6
$ws.Range("K86").Value = 1685.4445
$ws.Range("L86").Value = 11661.6
$ws.Range("M86").Value = -562.4445000000001
$ws.Range("N86").Value = -13907.6
$ws.Range("H89").Value = 6936.0527
$ws.Range("I89").Value = 1685.4445
$ws.Range("J89").Value = 11661.6
$ws.Range("K89").Value = 8427.2225
$ws.Range("L89").Value = 58308
$ws.Range("M89").Value = -2811.2225
$ws.Range("N89").Value = -69540
$ws.Range("H96").Value = 27778790
$ws.Range("I96").Value = 41667636
$ws.Range("J96").Value = 1096.6666
$ws.Range("K96").Value = 125002908
$ws.Range("L96").Value = 3289.9998
$ws.Range("M96").Value = -125001535
$ws.Range("N96").Value = -6035.9998
$ws.Range("H112").Value = 3704768.2
$ws.Range("J112").Value = 3969351.8
$ws.Range("L112").Value = 11908055.4
$ws.Range("N112").Value = -11910271.4
$ws.Range("H137").Value = 49485.617
$ws.Range("I137").Value = 1726.7273
$ws.Range("K137").Value = 5180.1819
$ws.Range("M137").Value = -2630.1819

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3474222.2
$ws.Range("I63").Value = 2250
$ws.Range("J63").Value = 31250000
$ws.Range("K63").Value = 2250
$ws.Range("L63").Value = 31250000
$ws.Range("M63").Value = -1564
$ws.Range("N63").Value = -31251372
$ws.Range("H66").Value = 3474222.2
$ws.Range("I66").Value = 2250
$ws.Range("J66").Value = 31250000
$ws.Range("K66").Value = 11250
$ws.Range("L66").Value = 156250000
$ws.Range("M66").Value = -7818
$ws.Range("N66").Value = -156256864
$ws.Range("H102").Value = 1266.3334
$ws.Range("I102").Value = 1083.1666
$ws.Range("J102").Value = 1999
$ws.Range("K102").Value = 1083.1666
$ws.Range("L102").Value = 1999
$ws.Range("M102").Value = 538.8334
$ws.Range("N102").Value = -5243
$ws.Range("H122").Value = 2445.0557
$ws.Range("I122").Value = 2469.4375
$ws.Range("J122").Value = 2250
$ws.Range("K122").Value = 7408.3125
$ws.Range("L122").Value = 6750
$ws.Range("M122").Value = -4958.3125
$ws.Range("N122").Value = -11650

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2514.7368
$ws.Range("I20").Value = 2267.6924
$ws.Range("K20").Value = 2267.6924
$ws.Range("M20").Value = -2020.6924

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 15600.409
$ws.Range("I31").Value = 21227.334
$ws.Range("J31").Value = 3542.7144
$ws.Range("K31").Value = 21227.334
$ws.Range("L31").Value = 3542.7144
$ws.Range("M31").Value = -20932.334
$ws.Range("N31").Value = -4132.7144
$ws.Range("H34").Value = 15600.409
$ws.Range("I34").Value = 21227.334
$ws.Range("J34").Value = 3542.7144
$ws.Range("K34").Value = 21227.334
$ws.Range("L34").Value = 3542.7144
$ws.Range("M34").Value = -21025.334
$ws.Range("N34").Value = -3946.7144
$ws.Range("H41").Value = 2750
$ws.Range("I41").Value = 2750
$ws.Range("K41").Value = 2750
$ws.Range("M41").Value = -2322
$ws.Range("H58").Value = 24975.047
$ws.Range("I58").Value = 1373.5
$ws.Range("K58").Value = 1373.5
$ws.Range("M58").Value = -1170.5
$ws.Range("H106").Value = 24671
$ws.Range("J106").Value = 24671
$ws.Range("L106").Value = 24671
$ws.Range("N106").Value = -27195
$ws.Range("H136").Value = 24975.047
$ws.Range("I136").Value = 1373.5
$ws.Range("K136").Value = 4120.5
$ws.Range("M136").Value = -1570.5

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 771.58
$ws.Range("J131").Value = 783.7553
$ws.Range("L131").Value = 2351.2659
$ws.Range("N131").Value = -12431.2659

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 3999.3333
$ws.Range("I43").Value = 3499
$ws.Range("J43").Value = 5000
$ws.Range("K43").Value = 3499
$ws.Range("L43").Value = 5000
$ws.Range("M43").Value = -3348
$ws.Range("N43").Value = -5302
$ws.Range("H70").Value = 3686811
$ws.Range("I70").Value = 18871.428
$ws.Range("J70").Value = 6254369
$ws.Range("K70").Value = 18871.428
$ws.Range("L70").Value = 6254369
$ws.Range("M70").Value = -18601.428
$ws.Range("N70").Value = -6254909
$ws.Range("H73").Value = 3686811
$ws.Range("I73").Value = 18871.428
$ws.Range("J73").Value = 6254369
$ws.Range("K73").Value = 18871.428
$ws.Range("L73").Value = 6254369
$ws.Range("M73").Value = -17935.428
$ws.Range("N73").Value = -6256241
$ws.Range("H80").Value = 3747.4666
$ws.Range("I80").Value = 3216.6667
$ws.Range("J80").Value = 4101.3335
$ws.Range("K80").Value = 3216.6667
$ws.Range("L80").Value = 4101.3335
$ws.Range("M80").Value = -2218.6667
$ws.Range("N80").Value = -6097.3335
$ws.Range("H83").Value = 3747.4666
$ws.Range("I83").Value = 3216.6667
$ws.Range("J83").Value = 4101.3335
$ws.Range("K83").Value = 16083.3335
$ws.Range("L83").Value = 20506.6675
$ws.Range("M83").Value = -11091.3335
$ws.Range("N83").Value = -30490.6675
$ws.Range("H122").Value = 2464.2666
$ws.Range("I122").Value = 2083.3333
$ws.Range("J122").Value = 2718.2222
$ws.Range("K122").Value = 6249.999899999999
$ws.Range("L122").Value = 8154.6666
$ws.Range("M122").Value = -3799.999899999999
$ws.Range("N122").Value = -13054.6666

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5706.1333
$ws.Range("I7").Value = 6141.9165
$ws.Range("K7").Value = 6141.9165
$ws.Range("M7").Value = -6029.9165
$ws.Range("H22").Value = 5268.1665
$ws.Range("I22").Value = 10001
$ws.Range("J22").Value = 4321.6
$ws.Range("K22").Value = 10001
$ws.Range("L22").Value = 4321.6
$ws.Range("M22").Value = -9706
$ws.Range("N22").Value = -4911.6
$ws.Range("H27").Value = 5268.1665
$ws.Range("I27").Value = 10001
$ws.Range("J27").Value = 4321.6
$ws.Range("K27").Value = 10001
$ws.Range("L27").Value = 4321.6
$ws.Range("M27").Value = -9894
$ws.Range("N27").Value = -4535.6
$ws.Range("H122").Value = 2805106.8
$ws.Range("I122").Value = 3924409.5
$ws.Range("J122").Value = 6850
$ws.Range("K122").Value = 11773228.5
$ws.Range("L122").Value = 20550
$ws.Range("M122").Value = -11770778.5
$ws.Range("N122").Value = -25450
$ws.Range("H125").Value = 39500
$ws.Range("J125").Value = 39500
$ws.Range("L125").Value = 39500
$ws.Range("N125").Value = -49340
$ws.Range("H126").Value = 5706.1333
$ws.Range("I126").Value = 6141.9165
$ws.Range("K126").Value = 18425.7495
$ws.Range("M126").Value = -15955.7495

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H44").Value = 17000
$ws.Range("J44").Value = 17000
$ws.Range("L44").Value = 17000
$ws.Range("N44").Value = -18108
$ws.Range("H107").Value = 1894550
$ws.Range("I107").Value = 640.6875
$ws.Range("K107").Value = 1922.0625
$ws.Range("M107").Value = -2.0625
$ws.Range("H132").Value = 2582.8333
$ws.Range("J132").Value = 2899.4
$ws.Range("L132").Value = 8698.200000000001
$ws.Range("N132").Value = -13758.2
$ws.Range("H136").Value = 31251986
$ws.Range("I136").Value = 55557490
$ws.Range("J136").Value = 2057.1428
$ws.Range("K136").Value = 166672470
$ws.Range("L136").Value = 6171.428400000001
$ws.Range("M136").Value = -166669920
$ws.Range("N136").Value = -11271.4284
